# Update the 6 worksheets: fill in B100/B101 with real values and append
# a new row 102 for date 45961 (2025-10-31) with a placeholder value of 0.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> (B100 value, B101 value)
$updates = @{
    "진양산업"   = @{ B100 = 3221; B101 = 3137 }
    "넥스트아이" = @{ B100 = 1188; B101 = 1174 }
    "삼보산업"   = @{ B100 = 1334; B101 = 1305 }
    "YBM넷"      = @{ B100 = 1912; B101 = 1900 }
    "NE능률"     = @{ B100 = 805;  B101 = 768  }
    "위즈코프"   = @{ B100 = 1540; B101 = 1526 }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]

        $ws.Range("B100").Value = $vals.B100
        $ws.Range("B101").Value = $vals.B101

        $ws.Range("A102").Value = 45961
        $ws.Range("A102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Range("B102").Value = 0
    }
}
